# Add a second "login" style entry (username/password test data) below the
# existing rows, mirroring the existing selenium / amazon-hyperlink layout:
#   C7 -> "aashi07"   (plain shared string, like A1)
#   H8 -> "ashi@123"  (shared string styled + hyperlinked, like E2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New plain text cell.
$ws.Range("C7").Value = "aashi07"

# New hyperlinked cell (value first, then turn it into a hyperlink).
$ws.Range("H8").Value = "ashi@123"
$null = $ws.Hyperlinks.Add($ws.Range("H8"), "ashi@123")

# Hyperlinks.Add() stamps its own style on the cell; force it back onto the
# workbook's shared "Hyperlink" cell style so H8 matches E2's formatting.
$ws.Range("H8").Style = "Hyperlink"

# Reflect the new active cell/selection recorded in the saved sheet view.
$null = $ws.Range("K6").Select()
